# Refresh the crypto price table with the latest scraped Price (col D) and
# Volume(1h) (col E) figures for each coin row.
#
# This run also has rows 32/33 swap identity: ImmutableX drops below
# Binance-Peg BSC-USD in the ranking, so for that pair the Coin name (B),
# Link (C), Price (D) and Volume (E) are all rewritten together.
#
# $updates maps each A1 address to its new value. Price-column values that
# look like plain numbers ("1.00", "682.02", ...) are prefixed with a
# leading "'" so Excel stores them as literal text, consistent with every
# other cell in this text-only worksheet, instead of silently coercing them
# to numerics. Because that apostrophe trick leaves a "quote prefix" number
# style attached to the cell, we immediately reset the range's style back
# to "Normal" afterwards so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "69.282.79"
    "E2" = "  -0.09%  "
    "D3" = "3.680.20"
    "E3" = "  -0.04%  "
    "D4" = "'1.00"
    "E4" = "  -0.01%  "
    "D5" = "'682.02"
    "E5" = "  -0.15%  "
    "D6" = "'157.97"
    "E6" = "  -2.84%  "
    "E7" = "  +0.02%  "
    "E8" = "  -1.27%  "
    "E9" = "  -1.60%  "
    "D10" = "'6.99"
    "E10" = "  -3.43%  "
    "E11" = "  -3.26%  "
    "E12" = "  -2.24%  "
    "D13" = "4.300.30"
    "E13" = "  -0.08%  "
    "D14" = "'32.21"
    "E14" = "  -4.09%  "
    "D15" = "3.679.34"
    "E15" = "  -0.20%  "
    "D16" = "69.314.66"
    "E16" = "  -0.07%  "
    "E17" = "  +1.72%  "
    "D18" = "'15.95"
    "E18" = "  -2.31%  "
    "E19" = "  -4.16%  "
    "D20" = "'470.45"
    "E20" = "  -1.95%  "
    "D21" = "'9.99"
    "E21" = "  +1.79%  "
    "D22" = "'0.649"
    "E22" = "  -2.67%  "
    "D23" = "'80.00"
    "E23" = "  -0.09%  "
    "D24" = "3.824.24"
    "E24" = "  -0.11%  "
    "E25" = "  -0.03%  "
    "D26" = "'0.0000122"
    "E26" = "  -5.07%  "
    "D27" = "'10.89"
    "E27" = "  -5.42%  "
    "D28" = "'9.12"
    "E28" = "  -5.06%  "
    "D29" = "'2.71"
    "E29" = "  -1.71%  "
    "E30" = "  -5.01%  "
    "E31" = "  -4.68%  "
    "B32" = "Binance-PegBSC-USD"
    "C32" = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
    "D32" = "'1.00"
    "E32" = "  +0.06%  "
    "B33" = "ImmutableX"
    "C33" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D33" = "'1.98"
    "E33" = "  -6.11%  "
    "D34" = "'26.84"
    "E34" = "  -1.06%  "
    "D35" = "3.659.41"
    "E35" = "  +0.40%  "
    "E36" = "  -4.44%  "
    "D37" = "'8.19"
    "E37" = "  -4.82%  "
    "D38" = "'6.03"
    "E38" = "  -1.71%  "
    "E40" = "  +2.46%  "
    "E41" = "  -4.38%  "
    "E42" = "  +0.02%  "
    "D43" = "'168.00"
    "E43" = "  +8.82%  "
    "D44" = "'0.940"
    "E44" = "  -1.98%  "
    "D45" = "'47.61"
    "E45" = "  -1.27%  "
    "D46" = "'2.71"
    "E46" = "  -4.93%  "
    "E47" = "  -1.47%  "
    "E48" = "  +2.07%  "
    "E49" = "  -5.02%  "
    "D50" = "'7.75"
    "E50" = "  -4.44%  "
    "D51" = "'27.03"
    "E51" = "  -3.15%  "
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $range = $ws.Range($addr)
    $range.Value = $newValue
    if ($newValue.Substring(0, 1) -eq [char]39) {
        $range.Style = "Normal"
    }
}
